$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.028185895818088
$ws.Range("D2").Value = 1.033314673539573
$ws.Range("E2").Value = 1.028181577463523
$ws.Range("F2").Value = 1.026732192368816
$ws.Range("I2").Value = 1.035243628196483
$ws.Range("J2").Value = 1.03333961427441
$ws.Range("K2").Value = 1.036117378245251
$ws.Range("L2").Value = 1.030999140757189
$ws.Range("M2").Value = 1.029553980392997

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.029130216393325
$ws.Range("D3").Value = 1.034040460955386
$ws.Range("E3").Value = 1.028983265508087
$ws.Range("F3").Value = 1.028309723001516
$ws.Range("I3").Value = 1.035506269172334
$ws.Range("J3").Value = 1.033924463486449
$ws.Range("K3").Value = 1.03665234238061
$ws.Range("L3").Value = 1.031608727094165
$ws.Range("M3").Value = 1.030937004167045

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.029741101974341
$ws.Range("D4").Value = 1.034509842357398
$ws.Range("E4").Value = 1.029502257531263
$ws.Range("F4").Value = 1.029330308996242
$ws.Range("I4").Value = 1.035674727158372
$ws.Range("J4").Value = 1.034302147432106
$ws.Range("K4").Value = 1.036997581361546
$ws.Range("L4").Value = 1.032002766181187
$ws.Range("M4").Value = 1.031831258504808

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.029997881914666
$ws.Range("D5").Value = 1.034707109708846
$ws.Range("E5").Value = 1.029720500419358
$ws.Range("F5").Value = 1.029759325292152
$ws.Range("I5").Value = 1.035745190567728
$ws.Range("J5").Value = 1.034460745724082
$ws.Range("K5").Value = 1.037142499865465
$ws.Range("L5").Value = 1.032168323436851
$ws.Range("M5").Value = 1.032207051008634

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.030040994272115
$ws.Range("D6").Value = 1.034740228201565
$ws.Range("E6").Value = 1.029757147790179
$ws.Range("F6").Value = 1.029831356943114
$ws.Range("I6").Value = 1.035757000788077
$ws.Range("J6").Value = 1.034487364521461
$ws.Range("K6").Value = 1.037166819419551
$ws.Range("L6").Value = 1.032196115551997
$ws.Range("M6").Value = 1.032270139444748

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.02974453322301
$ws.Range("D7").Value = 1.03451247849195
$ws.Range("E7").Value = 1.029505173473401
$ws.Range("F7").Value = 1.029336041675722
$ws.Range("I7").Value = 1.035675670094067
$ws.Range("J7").Value = 1.03430426733752
$ws.Range("K7").Value = 1.036999518634075
$ws.Range("L7").Value = 1.032004978745682
$ws.Range("M7").Value = 1.031836280454168

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.028505065213616
$ws.Range("D8").Value = 1.033560008674356
$ws.Range("E8").Value = 1.028452460518836
$ws.Range("F8").Value = 1.027265367172565
$ws.Range("I8").Value = 1.035332697087422
$ws.Range("J8").Value = 1.033537422718713
$ws.Range("K8").Value = 1.036298362049544
$ws.Range("L8").Value = 1.031205236943965
$ws.Range("M8").Value = 1.030021518092769

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.02631977997897
$ws.Range("D9").Value = 1.031879723389283
$ws.Range("E9").Value = 1.026599344785222
$ws.Range("F9").Value = 1.0236149366482
$ws.Range("I9").Value = 1.034716937151585
$ws.Range("J9").Value = 1.032180373729709
$ws.Range("K9").Value = 1.035055800882851
$ws.Range("L9").Value = 1.029792896781966
$ws.Range("M9").Value = 1.026818452567219

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.024862106045968
$ws.Range("D10").Value = 1.030758265516502
$ws.Range("E10").Value = 1.025365225911733
$ws.Range("F10").Value = 1.021179871528654
$ws.Range("I10").Value = 1.034298766789389
$ws.Range("J10").Value = 1.031271782480321
$ws.Range("K10").Value = 1.034222697327133
$ws.Range("L10").Value = 1.02884925661699
$ws.Range("M10").Value = 1.024679269998502

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.024230717670701
$ws.Range("D11").Value = 1.030272364074311
$ws.Range("E11").Value = 1.024831147435685
$ws.Range("F11").Value = 1.020125041885859
$ws.Range("I11").Value = 1.034115875981356
$ws.Range("J11").Value = 1.030877425873767
$ws.Range("K11").Value = 1.033860831616492
$ws.Range("L11").Value = 1.028440155700228
$ws.Range("M11").Value = 1.023752011775769

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.023996160373102
$ws.Range("D12").Value = 1.030091833289079
$ws.Range("E12").Value = 1.024632812542575
$ws.Range("F12").Value = 1.019733159843322
$ws.Range("I12").Value = 1.034047668579912
$ws.Range("J12").Value = 1.030730804068822
$ws.Range("K12").Value = 1.033726249178814
$ws.Range("L12").Value = 1.028288122337454
$ws.Range("M12").Value = 1.023407434091709

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.024046475146766
$ws.Range("D13").Value = 1.030130559805678
$ws.Range("E13").Value = 1.024675353984097
$ws.Range("F13").Value = 1.019817223208011
$ws.Range("I13").Value = 1.034062311682797
$ws.Range("J13").Value = 1.030762261294512
$ws.Range("K13").Value = 1.033755125237187
$ws.Range("L13").Value = 1.028320737409665
$ws.Range("M13").Value = 1.023481354176991

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.024211329738275
$ws.Range("D14").Value = 1.030257442267768
$ws.Range("E14").Value = 1.024814752080172
$ws.Range("F14").Value = 1.020092650288218
$ws.Range("I14").Value = 1.034110243518676
$ws.Range("J14").Value = 1.030865308932007
$ws.Range("K14").Value = 1.033849710446977
$ws.Range("L14").Value = 1.028427590107409
$ws.Range("M14").Value = 1.023723532025188

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.024312897895187
$ws.Range("D15").Value = 1.030335612734222
$ws.Range("E15").Value = 1.02490064589718
$ws.Range("F15").Value = 1.020262340347859
$ws.Range("I15").Value = 1.034139739644337
$ws.Range("J15").Value = 1.030928781392465
$ws.Range("K15").Value = 1.033907965063816
$ws.Range("L15").Value = 1.028493415628645
$ws.Range("M15").Value = 1.023872725384396

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.024904004822816
$ws.Range("D16").Value = 1.03079050682739
$ws.Range("E16").Value = 1.025400677335894
$ws.Range("F16").Value = 1.02124986746034
$ws.Range("I16").Value = 1.034310866287268
$ws.Range("J16").Value = 1.031297934989692
$ws.Range("K16").Value = 1.034246689416343
$ws.Range("L16").Value = 1.028876396780085
$ws.Range("M16").Value = 1.024740787975994

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.025274735022631
$ws.Range("D17").Value = 1.031075768908843
$ws.Range("E17").Value = 1.025714415072168
$ws.Range("F17").Value = 1.021869197193426
$ws.Range("I17").Value = 1.034417722073953
$ws.Range("J17").Value = 1.031529245867326
$ws.Range("K17").Value = 1.034458860570037
$ws.Range("L17").Value = 1.029116497156257
$ws.Range("M17").Value = 1.025285034801428

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.025490955846199
$ws.Range("D18").Value = 1.031242128232068
$ws.Range("E18").Value = 1.025897442228874
$ws.Range("F18").Value = 1.022230400673474
$ws.Range("I18").Value = 1.034479873624973
$ws.Range("J18").Value = 1.031664075826004
$ws.Range("K18").Value = 1.03458250774785
$ws.Range("L18").Value = 1.029256495557701
$ws.Range("M18").Value = 1.025602390785901

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.025564678214112
$ws.Range("D19").Value = 1.031298847501256
$ws.Range("E19").Value = 1.025959854716748
$ws.Range("F19").Value = 1.022353554961805
$ws.Range("I19").Value = 1.034501035913339
$ws.Range("J19").Value = 1.031710034144033
$ws.Range("K19").Value = 1.034624649813668
$ws.Range("L19").Value = 1.02930422324368
$ws.Range("M19").Value = 1.025710585209229

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.025234961235991
$ws.Range("D20").Value = 1.031045166008725
$ws.Range("E20").Value = 1.025680750956715
$ws.Range("F20").Value = 1.021802753270503
$ws.Range("I20").Value = 1.034406275625564
$ws.Range("J20").Value = 1.031504437686521
$ws.Range("K20").Value = 1.03443610785424
$ws.Range("L20").Value = 1.029090741633043
$ws.Range("M20").Value = 1.025226652029547

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.024162785045563
$ws.Range("D21").Value = 1.030220079796126
$ws.Range("E21").Value = 1.024773701554925
$ws.Range("F21").Value = 1.020011545901315
$ws.Range("I21").Value = 1.03409613634182
$ws.Range("J21").Value = 1.030834967844609
$ws.Range("K21").Value = 1.033821862137905
$ws.Range("L21").Value = 1.028396126724805
$ws.Range("M21").Value = 1.023652220946443

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.023488482521963
$ws.Range("D22").Value = 1.029701052997205
$ws.Range("E22").Value = 1.02420366764982
$ws.Range("F22").Value = 1.018884926611476
$ws.Range("I22").Value = 1.033899556359579
$ws.Range("J22").Value = 1.030413234216386
$ws.Range("K22").Value = 1.03343468128417
$ws.Range("L22").Value = 1.02795896035657
$ws.Range("M22").Value = 1.022661427328569

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.023845960578084
$ws.Range("D23").Value = 1.029976223748947
$ws.Range("E23").Value = 1.02450582835166
$ws.Range("F23").Value = 1.019482210498729
$ws.Range("I23").Value = 1.034003917233916
$ws.Range("J23").Value = 1.030636880232874
$ws.Range("K23").Value = 1.033640026152118
$ws.Range("L23").Value = 1.02819075175648
$ws.Range("M23").Value = 1.02318675171167

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.02525293335914
$ws.Range("D24").Value = 1.031058994232184
$ws.Range("E24").Value = 1.02569596223096
$ws.Range("F24").Value = 1.021832776544583
$ws.Range("I24").Value = 1.034411448325418
$ws.Range("J24").Value = 1.031515647712729
$ws.Range("K24").Value = 1.034446389162739
$ws.Range("L24").Value = 1.029102379593777
$ws.Range("M24").Value = 1.025253032980122

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.026884871056713
$ws.Range("D25").Value = 1.032314342108266
$ws.Range("E25").Value = 1.027078193018232
$ws.Range("F25").Value = 1.024558889811962
$ws.Range("I25").Value = 1.034877476583388
$ws.Range("J25").Value = 1.032531888495359
$ws.Range("K25").Value = 1.035377866268381
$ws.Range("L25").Value = 1.030158387123598
$ws.Range("M25").Value = 1.027647170896795
